# Generate Report for Handback
# This script updates the localization-status workbook to reflect a
# completed handback from the de-de / zh-cn localization pipelines.

$wb = $excel.ActiveWorkbook

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f270790db99518d8c5b55d52635601fabe283b/e2e/5b671294-65cc-40f4-bdbc-50b430ae538b.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f270790db99518d8c5b55d52635601fabe283b/e2e/a965fe83-e1bc-49b9-bc23-10cd6eee590d.md"
$disp1 = "5b671294-65cc-40f4-bdbc-50b430ae538b.md"
$disp2 = "a965fe83-e1bc-49b9-bc23-10cd6eee590d.md"

# ---------------------------------------------------------------------
# Overview sheet: the Status text used across the workbook changes from
# "Ready for handoff" to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.09
$wsOverview.Columns.Item(6).ColumnWidth = 29.09

# ---------------------------------------------------------------------
# zh-cn sheet: fill in the "Latest Target File" / "Latest Handback File"
# columns and wire up hyperlinks for the handback markdown files.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("J2").Value = "5b671294-65cc-40f4-bdbc-50b430ae538b.15125abe94f394fe1e89425320597259bbe931a6.zh-cn.xlf"
$wsZh.Range("J3").Value = "a965fe83-e1bc-49b9-bc23-10cd6eee590d.690fa03af6d97ad3909657bbbd05d96f8eec8091.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-30 20:53:53"
$wsZh.Range("K3").Value = "2016-08-30 20:53:53"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlMd1, $null, $null, $disp1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlMd1, $null, $null, $disp1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlMd2, $null, $null, $disp2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlMd2, $null, $null, $disp2)

$wsZh.Columns.Item(3).ColumnWidth = 29.09
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, but the handback datetime
# is a freshly introduced value (was the zero-date placeholder before).
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("J2").Value = "5b671294-65cc-40f4-bdbc-50b430ae538b.15125abe94f394fe1e89425320597259bbe931a6.de-de.xlf"
$wsDe.Range("J3").Value = "a965fe83-e1bc-49b9-bc23-10cd6eee590d.690fa03af6d97ad3909657bbbd05d96f8eec8091.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-30 20:54:01"
$wsDe.Range("K3").Value = "2016-08-30 20:54:01"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlMd1, $null, $null, $disp1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlMd1, $null, $null, $disp1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlMd2, $null, $null, $disp2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlMd2, $null, $null, $disp2)

$wsDe.Columns.Item(3).ColumnWidth = 29.09
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
